$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F3").Value = 0.8
$ws.Range("E2").Value = "Input how much money you want to put in this weight of stock"
$ws.Range("E2").Font.Bold = $true
$ws.Range("F7").Interior.Color = 15773696
$ws.Range("H15:H16").Select()
$ws.Application.ActiveWindow.RangeSelection.Activate()
